$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The indicator rows in the "MSME Participation on the Economy" block (rows
# 10-14) need their backing shared-string entries reordered so that
# "Enterprises density (per 1000 people)" / "33.9" sit right after "MSMEs"
# (i.e. before "Employment (% of total)") instead of after "180000". The
# on-screen layout/values do not change at all - only the order in which
# the strings were (re)created, which drives their position in the shared
# string table.
#
# Clear the whole block first so every label/value currently in it is
# dropped from the shared-string pool, then re-enter them in the desired
# order so they get re-interned in that order.
$ws.Range("A10:D14").ClearContents()

# Row 10: Enterprises density (per 1000 people) / 33.9
$ws.Range("A10").Value = "Enterprises density (per 1000 people)"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.9"
$ws.Range("D10").Style = "Normal"

# Row 11: Employment (% of total) / 67.7
$ws.Range("A11").Value = "Employment (% of total)"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "67.7"
$ws.Range("D11").Style = "Normal"

# Row 12: Enterprises (absolute #) / 180000
$ws.Range("A12").Value = "Enterprises (absolute #)"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "180000"
$ws.Range("D12").Style = "Normal"

# Row 13: Employment (absolute #) / 2100000
$ws.Range("A13").Value = "Employment (absolute #)"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2100000"
$ws.Range("D13").Style = "Normal"

# Row 14: Enterprises (% of total) / 99.3
$ws.Range("A14").Value = "Enterprises (% of total)"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "99.3"
$ws.Range("D14").Style = "Normal"

# Restore the label column's style (bold "title_" look) that ClearContents
# left untouched on A, but make sure A column keeps its original style too.
$ws.Range("A10:A14").Style = "title_"
